$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 824.75
$ws.Range("I21").Value = 649.5
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 649.5
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = -181.5
$ws.Range("N21").Value = -1936

$ws.Range("H23").Value = 824.75
$ws.Range("I23").Value = 649.5
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 649.5
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = -415.5
$ws.Range("N23").Value = -1468

$ws.Range("H29").Value = 4862.5
$ws.Range("I29").Value = 175
$ws.Range("J29").Value = 5800
$ws.Range("K29").Value = 525
$ws.Range("L29").Value = 17400
$ws.Range("M29").Value = -244
$ws.Range("N29").Value = -17962

$ws.Range("H38").Value = 1965.4
$ws.Range("I38").Value = 43.727272
$ws.Range("J38").Value = 7250
$ws.Range("K38").Value = 131.181816
$ws.Range("L38").Value = 21750
$ws.Range("M38").Value = 240.818184
$ws.Range("N38").Value = -22494

$ws.Range("H40").Value = 7777
$ws.Range("I40").Value = 5999.4
$ws.Range("K40").Value = 5999.4
$ws.Range("M40").Value = -5824.4

$ws.Range("H58").Value = 2167.8667
$ws.Range("I58").Value = 1112.3334
$ws.Range("J58").Value = 2871.5557
$ws.Range("K58").Value = 3337.0002
$ws.Range("L58").Value = 8614.667099999999
$ws.Range("M58").Value = -3187.0002
$ws.Range("N58").Value = -8914.667099999999

$ws.Range("H86").Value = 2929.1428
$ws.Range("J86").Value = 3000.6667
$ws.Range("L86").Value = 3000.6667
$ws.Range("N86").Value = -5246.6667

$ws.Range("H89").Value = 2929.1428
$ws.Range("J89").Value = 3000.6667
$ws.Range("L89").Value = 15003.3335
$ws.Range("N89").Value = -26235.3335

$ws.Range("H135").Value = 907.9091
$ws.Range("I135").Value = 743.55554
$ws.Range("K135").Value = 6691.99986
$ws.Range("M135").Value = -4156.99986

$ws.Range("H137").Value = 2284.9375
$ws.Range("I137").Value = 1090.4445
$ws.Range("K137").Value = 3271.3335
$ws.Range("M137").Value = -721.3335000000002


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H86").Value = 5907.4165
$ws.Range("I86").Value = 2148.3333
$ws.Range("K86").Value = 2148.3333
$ws.Range("M86").Value = -1025.3333

$ws.Range("H89").Value = 5907.4165
$ws.Range("I89").Value = 2148.3333
$ws.Range("K89").Value = 10741.6665
$ws.Range("M89").Value = -5125.666499999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2256.111
$ws.Range("J22").Value = 3721.75
$ws.Range("L22").Value = 3721.75
$ws.Range("N22").Value = -4421.75

$ws.Range("H58").Value = 2240.5
$ws.Range("I58").Value = 1617.7646
$ws.Range("K58").Value = 1617.7646
$ws.Range("M58").Value = -1414.7646

$ws.Range("H136").Value = 2240.5
$ws.Range("I136").Value = 1617.7646
$ws.Range("K136").Value = 4853.293799999999
$ws.Range("M136").Value = -2303.293799999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1750
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9920

$ws.Range("H104").Value = 11497.143
$ws.Range("J104").Value = 11996.923
$ws.Range("L104").Value = 35990.769
$ws.Range("N104").Value = -41232.769


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9000
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8730
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 9000
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8064
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 3239.5
$ws.Range("I80").Value = 3239.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3239.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2241.5
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3239.5
$ws.Range("I83").Value = 3239.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 16197.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11205.5
$ws.Range("N83").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2028.2858
$ws.Range("I16").Value = 2032.6666
$ws.Range("J16").Value = 2002
$ws.Range("K16").Value = 2032.6666
$ws.Range("L16").Value = 2002
$ws.Range("M16").Value = -1862.6666
$ws.Range("N16").Value = -2342

$ws.Range("H82").Value = 3259.9285
$ws.Range("I82").Value = 611.25
$ws.Range("K82").Value = 611.25
$ws.Range("M82").Value = -250.25

$ws.Range("H85").Value = 3259.9285
$ws.Range("I85").Value = 611.25
$ws.Range("K85").Value = 611.25
$ws.Range("M85").Value = 636.75

$ws.Range("H132").Value = 963.5
$ws.Range("I132").Value = 984.6667
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 2954.0001
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -424.0001000000002
$ws.Range("N132").Value = -7760


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 290
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H9").Value = 499.66666
$ws.Range("I9").Value = 599
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 599
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = -459
$ws.Range("N9").Value = -730

$ws.Range("H14").Value = 17500
$ws.Range("I14").Value = 17500
$ws.Range("K14").Value = 17500
$ws.Range("M14").Value = -17332

$ws.Range("H34").Value = 11458.667
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H126").Value = 4511
$ws.Range("I126").Value = 2465.7778
$ws.Range("J126").Value = 7140.5713
$ws.Range("K126").Value = 7397.3334
$ws.Range("L126").Value = 21421.7139
$ws.Range("M126").Value = -4927.3334
$ws.Range("N126").Value = -26361.7139

